$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.052.55"
$ws.Range("E2").Value = "  +2.73%  "

$ws.Range("D3").Value = "2.057.04"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.30"
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.36"
$ws.Range("E7").Value = "  +6.85%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +3.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0808"
$ws.Range("E10").Value = "  +3.47%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").Value = "2.361.81"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.65"
$ws.Range("E13").Value = "  +4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.71"
$ws.Range("E14").Value = "  +2.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.754"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +3.23%  "

$ws.Range("D17").Value = "2.063.74"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("D18").Value = "37.933.13"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.82"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.93"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("E25").Value = "  +4.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.32"
$ws.Range("E26").Value = "  +2.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.59"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("E28").Value = "  +8.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.06"
$ws.Range("E29").Value = "  +2.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.55"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +4.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0613"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.01"
$ws.Range("E35").Value = "  +9.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.99"
$ws.Range("E37").Value = "  +13.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +5.66%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.47"
$ws.Range("E40").Value = "  +4.09%  "

$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("D42").Value = "1.484.37"
$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("E43").Value = "  +3.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0938"
$ws.Range("E44").Value = "  +2.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.75"
$ws.Range("E45").Value = "  +3.54%  "

$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.11"
$ws.Range("E47").Value = "  +17.71%  "

$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.06"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("D51").Value = "2.247.10"
$ws.Range("E51").Value = "  +2.14%  "
